# chore: update Sheets via scheduled runner
#
# Applies the leve-profit recalculation updates captured in the commit's
# diff: per-row market-price / profit figures changed on several of the
# crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), including a
# handful of rows that previously had all-zero price columns and now gain
# real values (plus newly-populated Profit columns M/N), and one row whose
# LeveProfitNQ (M55 on ALC) is cleared entirely.
#
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ,
#          J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ,
#          M=LeveProfitNQ, N=LeveProfitHQ
#
# Pass $null for a column that must stay untouched. Pass the string
# "CLEAR" for a column whose cell must be removed entirely (blanked).

$wb = $excel.ActiveWorkbook

function Set-LeveRow($SheetName, $Row, $H, $I, $J, $K, $L, $M, $N) {
    $ws = $wb.Worksheets.Item($SheetName)

    if ($H -eq "CLEAR") { $ws.Cells.Item($Row, 8).ClearContents() }
    elseif ($null -ne $H) { $ws.Cells.Item($Row, 8).Value = $H }

    if ($I -eq "CLEAR") { $ws.Cells.Item($Row, 9).ClearContents() }
    elseif ($null -ne $I) { $ws.Cells.Item($Row, 9).Value = $I }

    if ($J -eq "CLEAR") { $ws.Cells.Item($Row, 10).ClearContents() }
    elseif ($null -ne $J) { $ws.Cells.Item($Row, 10).Value = $J }

    if ($K -eq "CLEAR") { $ws.Cells.Item($Row, 11).ClearContents() }
    elseif ($null -ne $K) { $ws.Cells.Item($Row, 11).Value = $K }

    if ($L -eq "CLEAR") { $ws.Cells.Item($Row, 12).ClearContents() }
    elseif ($null -ne $L) { $ws.Cells.Item($Row, 12).Value = $L }

    if ($M -eq "CLEAR") { $ws.Cells.Item($Row, 13).ClearContents() }
    elseif ($null -ne $M) { $ws.Cells.Item($Row, 13).Value = $M }

    if ($N -eq "CLEAR") { $ws.Cells.Item($Row, 14).ClearContents() }
    elseif ($null -ne $N) { $ws.Cells.Item($Row, 14).Value = $N }
}

# ---------------------------------------------------------------- ALC ----
Set-LeveRow "ALC" 55  96.333336   0           96.333336   0           96.333336   "CLEAR"     -524.333336
Set-LeveRow "ALC" 98  751.6       751.6       $null       751.6       $null       746.4       $null
Set-LeveRow "ALC" 121 1883.5834  $null        1883.5834   $null       5650.7502   $null       -9144.7502
Set-LeveRow "ALC" 122 751.6       751.6       $null       2254.8      $null       195.1999999999998 $null
Set-LeveRow "ALC" 132 1004.9259   1017.56     847         3052.68     2541        -522.6799999999998 -7601
Set-LeveRow "ALC" 135 1858        1636.2727   2467.75     14726.4543  22209.75    -12191.4543 -27279.75
Set-LeveRow "ALC" 137 1840.4073  $null        2950.4167   $null       8851.250100000001 $null -13951.2501

# ---------------------------------------------------------------- ARM ----
Set-LeveRow "ARM" 2   1999.5      1999        2000        1999        2000        -1886       -2226
Set-LeveRow "ARM" 32  2249.7778   1945.4783   $null       1945.4783   $null       -1658.4783  $null
Set-LeveRow "ARM" 45  1994.5      1994.5      $null       1994.5      $null       -1617.5     $null
Set-LeveRow "ARM" 61  5728.2856   5728.2856   $null       5728.2856   $null       -5516.2856  $null
Set-LeveRow "ARM" 74  605         566.7       $null       566.7       $null       307.3       $null
Set-LeveRow "ARM" 77  605         566.7       $null       2833.5      $null       1534.5      $null
Set-LeveRow "ARM" 116 1999.5      1999        2000        1999        2000        295         -6588
Set-LeveRow "ARM" 132 3246.1765   2671.818    $null       8015.454000000001 $null -5485.454000000001 $null
Set-LeveRow "ARM" 136 5728.2856   5728.2856   $null       17184.8568  $null       -14634.8568 $null

# ---------------------------------------------------------------- BSM ----
Set-LeveRow "BSM" 3   1999.5      1999        2000        1999        2000        -1885       -2228
Set-LeveRow "BSM" 28  99500       $null       99500       $null       99500       $null       -100088
Set-LeveRow "BSM" 96  33999       33999       $null       33999       $null       -31253      $null
Set-LeveRow "BSM" 107 1125        2000        250         2000        250         -80         -4090
Set-LeveRow "BSM" 134 2889.818    2660.842    4340        7982.526    13020       -5447.526   -18090

# ---------------------------------------------------------------- CRP ----
Set-LeveRow "CRP" 31  1096.381    1044.0667   1227.1666   1044.0667   1227.1666   -749.0667000000001 -1817.1666
Set-LeveRow "CRP" 34  1096.381    1044.0667   1227.1666   1044.0667   1227.1666   -842.0667000000001 -1631.1666
Set-LeveRow "CRP" 58  2149.647    2103.0667   $null       2103.0667   $null       -1900.0667  $null
Set-LeveRow "CRP" 92  18000       $null       18000       $null       18000       $null       -22992
Set-LeveRow "CRP" 94  2995.7144   3000        2992.5      3000        2992.5      -2549       -3894.5
Set-LeveRow "CRP" 132 3927.5925   3567.4348   $null       10702.3044  $null       -8172.304400000001 $null
Set-LeveRow "CRP" 134 3446.96     3446.8262   3448.5      10340.4786  10345.5     -7805.4786  -15415.5
Set-LeveRow "CRP" 136 2149.647    2103.0667   $null       6309.2001   $null       -3759.2001  $null

# ---------------------------------------------------------------- CUL ----
Set-LeveRow "CUL" 107 1509.375    10000       296.42856   30000       889.28568   -28080      -4729.28568
Set-LeveRow "CUL" 131 608.1539    608.1539    $null       1824.4617   $null       3215.5383   $null

# ---------------------------------------------------------------- GSM ----
Set-LeveRow "GSM" 132 2477.3572   2117.2727   3797.6667   6351.8181   11393.0001  -3821.8181  -16453.0001

# ---------------------------------------------------------------- LTW ----
Set-LeveRow "LTW" 107 4444        4444        $null       4444        $null       -2524       $null

# ---------------------------------------------------------------- WVR ----
Set-LeveRow "WVR" 132 2036.8966   1258.0588   $null       3774.1764   $null       -1244.1764  $null
Set-LeveRow "WVR" 136 7500        5000        10000       15000       30000       -12450      -35100
